$d = $word.ActiveDocument

# 1) AstEvaluator.java:192 -> 189
$r1 = $d.Content.Find.Execute('AstEvaluator.java:192', $false, $false, $false, $false, $false, $true, 1, $false, 'AstEvaluator.java:189', 2)

# 2) AstEvaluator.java:112 -> 109
$r2 = $d.Content.Find.Execute('AstEvaluator.java:112', $false, $false, $false, $false, $false, $true, 1, $false, 'AstEvaluator.java:109', 2)

# 3) AbstractTemplatesTestSuite.java:461 -> 480
$r3 = $d.Content.Find.Execute('AbstractTemplatesTestSuite.java:461', $false, $false, $false, $false, $false, $true, 1, $false, 'AbstractTemplatesTestSuite.java:480', 2)

# 4) AbstractTemplatesTestSuite.java:368 -> 389
$r4 = $d.Content.Find.Execute('AbstractTemplatesTestSuite.java:368', $false, $false, $false, $false, $false, $true, 1, $false, 'AbstractTemplatesTestSuite.java:389', 2)

# 5) Insert 16 duplicated stack-trace lines after the *second* RunAfters.java:27) occurrence
$anchorFind = 'RunAfters.java:27)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)'
$anchorReplace = 'RunAfters.java:27)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.junit.runners.Suite.runChild(Suite.java:128)
	at org.junit.runners.Suite.runChild(Suite.java:27)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.junit.runners.Suite.runChild(Suite.java:128)
	at org.junit.runners.Suite.runChild(Suite.java:27)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)'
$r5 = $d.Content.Find.Execute($anchorFind, $false, $false, $false, $false, $false, $true, 1, $false, $anchorReplace, 2)

Write-Output "r1=$r1 r2=$r2 r3=$r3 r4=$r4 r5=$r5"
